# Edit script: updates "Product backlog" sheet content and refreshes the
# active-sheet / selection state on both "Product backlog" and
# "Sprint backlog" sheets.

$wb = $excel.ActiveWorkbook

$wsProduct = $wb.Worksheets.Item("Product backlog")
$wsSprint  = $wb.Worksheets.Item("Sprint backlog")

# --- Product backlog: header row -------------------------------------
# C3 ("Points" label) is cleared out.
$wsProduct.Range("C3").ClearContents()

# --- Product backlog: body rows ---------------------------------------
# Row 5 now shows "Thong bao trang thai ve" with 3 points.
$wsProduct.Range("B5").Value = "Thông báo trạng thái vé"
$wsProduct.Range("C5").Value = 3

# Row 7 now shows "Ve usecase xac dinh usecase" with 2 points (and the
# points cell loses its border formatting).
$wsProduct.Range("B7").Value = "Vẽ usecase xác định usecase"
$wsProduct.Range("C7").ClearFormats()
$wsProduct.Range("C7").Value = 2

# Row 8 now shows "Ve activity" with 2 points (points cell also loses its
# border formatting).
$wsProduct.Range("B8").Value = "Vẽ activity"
$wsProduct.Range("C8").ClearFormats()
$wsProduct.Range("C8").Value = 2

# Rows 9-14 are emptied out (contents only, formatting/styles stay).
$wsProduct.Range("A9:C14").ClearContents()

# --- Sheet selection / active tab -------------------------------------
# Sprint backlog keeps a remembered selection of C9 even though it is no
# longer the active sheet.
$wsSprint.Range("C9").Select()

# Product backlog becomes the active sheet, with D12 selected.
$wsProduct.Range("D12").Select()
